$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1256.619
$ws.Range("I17").Value = 1215
$ws.Range("J17").Value = 1261
$ws.Range("K17").Value = 3645
$ws.Range("L17").Value = 3783
$ws.Range("M17").Value = -3477
$ws.Range("N17").Value = -4119

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 5684.4443
$ws.Range("J121").Value = 5000.7144
$ws.Range("L121").Value = 15002.1432
$ws.Range("N121").Value = -18496.1432

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4853.5903
$ws.Range("I132").Value = 2633.0364
$ws.Range("J132").Value = 25208.666
$ws.Range("K132").Value = 7899.1092
$ws.Range("L132").Value = 75625.99800000001
$ws.Range("M132").Value = -5369.1092
$ws.Range("N132").Value = -80685.99800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 7739.25
$ws.Range("I14").Value = 6983.3335
$ws.Range("K14").Value = 6983.3335
$ws.Range("M14").Value = -6808.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 5885.8
$ws.Range("I21").Value = 4857.25
$ws.Range("K21").Value = 4857.25
$ws.Range("M21").Value = -4483.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7093.4575
$ws.Range("I32").Value = 6824.3965
$ws.Range("K32").Value = 6824.3965
$ws.Range("M32").Value = -6537.3965

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7247.161
$ws.Range("I45").Value = 8325.182000000001
$ws.Range("J45").Value = 4612
$ws.Range("K45").Value = 8325.182000000001
$ws.Range("L45").Value = 4612
$ws.Range("M45").Value = -7948.182000000001
$ws.Range("N45").Value = -5366

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3379.2942
$ws.Range("I61").Value = 3278.0625
$ws.Range("K61").Value = 3278.0625
$ws.Range("M61").Value = -3066.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3271.3572
$ws.Range("I63").Value = 1959.8
$ws.Range("K63").Value = 1959.8
$ws.Range("M63").Value = -1273.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3271.3572
$ws.Range("I66").Value = 1959.8
$ws.Range("K66").Value = 9799
$ws.Range("M66").Value = -6367

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 86496
$ws.Range("J122").Value = 86496
$ws.Range("L122").Value = 259488
$ws.Range("N122").Value = -264388

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2260.2727
$ws.Range("I132").Value = 2309.0625
$ws.Range("J132").Value = 699
$ws.Range("K132").Value = 6927.1875
$ws.Range("L132").Value = 2097
$ws.Range("M132").Value = -4397.1875
$ws.Range("N132").Value = -7157

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3379.2942
$ws.Range("I136").Value = 3278.0625
$ws.Range("K136").Value = 9834.1875
$ws.Range("M136").Value = -7284.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1936.2069
$ws.Range("I20").Value = 1561.7727
$ws.Range("K20").Value = 1561.7727
$ws.Range("M20").Value = -1314.7727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4866.8237
$ws.Range("I105").Value = 2662.4
$ws.Range("K105").Value = 2662.4
$ws.Range("M105").Value = -915.4000000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2525.5
$ws.Range("I107").Value = 2052.3333
$ws.Range("J107").Value = 3945
$ws.Range("K107").Value = 2052.3333
$ws.Range("L107").Value = 3945
$ws.Range("M107").Value = -132.3332999999998
$ws.Range("N107").Value = -7785

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 56.75
$ws.Range("I7").Value = 26.75
$ws.Range("J7").Value = 86.75
$ws.Range("K7").Value = 26.75
$ws.Range("L7").Value = 86.75
$ws.Range("M7").Value = 86.25
$ws.Range("N7").Value = -312.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2712.8572
$ws.Range("J31").Value = 7500
$ws.Range("L31").Value = 7500
$ws.Range("N31").Value = -8090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2712.8572
$ws.Range("J34").Value = 7500
$ws.Range("L34").Value = 7500
$ws.Range("N34").Value = -7904

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 29000
$ws.Range("I55").Value = 29000
$ws.Range("K55").Value = 29000
$ws.Range("M55").Value = -28685

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3332.7
$ws.Range("I58").Value = 3807.2856
$ws.Range("J58").Value = 2225.3333
$ws.Range("K58").Value = 3807.2856
$ws.Range("L58").Value = 2225.3333
$ws.Range("M58").Value = -3604.2856
$ws.Range("N58").Value = -2631.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10901.907
$ws.Range("I99").Value = 7796.3105
$ws.Range("J99").Value = 17334.928
$ws.Range("K99").Value = 7796.3105
$ws.Range("L99").Value = 17334.928
$ws.Range("M99").Value = -6298.3105
$ws.Range("N99").Value = -20330.928

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10901.907
$ws.Range("I126").Value = 7796.3105
$ws.Range("J126").Value = 17334.928
$ws.Range("K126").Value = 23388.9315
$ws.Range("L126").Value = 52004.784
$ws.Range("M126").Value = -20918.9315
$ws.Range("N126").Value = -56944.784

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5590.048
$ws.Range("I132").Value = 1786.129
$ws.Range("J132").Value = 16310.182
$ws.Range("K132").Value = 5358.387
$ws.Range("L132").Value = 48930.546
$ws.Range("M132").Value = -2828.387
$ws.Range("N132").Value = -53990.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2767.6511
$ws.Range("I134").Value = 2654.3845
$ws.Range("J134").Value = 3872
$ws.Range("K134").Value = 7963.1535
$ws.Range("L134").Value = 11616
$ws.Range("M134").Value = -5428.1535
$ws.Range("N134").Value = -16686

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3332.7
$ws.Range("I136").Value = 3807.2856
$ws.Range("J136").Value = 2225.3333
$ws.Range("K136").Value = 11421.8568
$ws.Range("L136").Value = 6675.999899999999
$ws.Range("M136").Value = -8871.856800000001
$ws.Range("N136").Value = -11775.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 112499.336
$ws.Range("J37").Value = 112499.336
$ws.Range("L37").Value = 337498.008
$ws.Range("N37").Value = -337722.008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 527.1429000000001
$ws.Range("I86").Value = 351.2857
$ws.Range("K86").Value = 1053.8571
$ws.Range("M86").Value = 132.1428999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 527.1429000000001
$ws.Range("I89").Value = 351.2857
$ws.Range("K89").Value = 3161.5713
$ws.Range("M89").Value = 2766.4287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 608.8570999999999
$ws.Range("I92").Value = 738.1667
$ws.Range("J92").Value = 573.5909
$ws.Range("K92").Value = 2214.5001
$ws.Range("L92").Value = 1720.7727
$ws.Range("M92").Value = -966.5001000000002
$ws.Range("N92").Value = -4216.7727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 5136.5
$ws.Range("I134").Value = 4791.1
$ws.Range("K134").Value = 14373.3
$ws.Range("M134").Value = -9303.300000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4064.25
$ws.Range("I70").Value = 4085.6667
$ws.Range("K70").Value = 4085.6667
$ws.Range("M70").Value = -3815.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4064.25
$ws.Range("I73").Value = 4085.6667
$ws.Range("K73").Value = 4085.6667
$ws.Range("M73").Value = -3149.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 17236.25
$ws.Range("J92").Value = 17236.25
$ws.Range("L92").Value = 17236.25
$ws.Range("N92").Value = -20980.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 6417.1177
$ws.Range("I107").Value = 626.25
$ws.Range("J107").Value = 11564.556
$ws.Range("K107").Value = 626.25
$ws.Range("L107").Value = 11564.556
$ws.Range("M107").Value = 1293.75
$ws.Range("N107").Value = -15404.556

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2420.72
$ws.Range("I132").Value = 1995.875
$ws.Range("J132").Value = 3176
$ws.Range("K132").Value = 5987.625
$ws.Range("L132").Value = 9528
$ws.Range("M132").Value = -3457.625
$ws.Range("N132").Value = -14588

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 16499.9
$ws.Range("J21").Value = 16499.9
$ws.Range("L21").Value = 16499.9
$ws.Range("N21").Value = -16847.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 12501456
$ws.Range("I93").Value = 12501456
$ws.Range("K93").Value = 12501456
$ws.Range("M93").Value = -12500208

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3422.4482
$ws.Range("I122").Value = 9836.6
$ws.Range("J122").Value = 2086.1667
$ws.Range("K122").Value = 29509.8
$ws.Range("L122").Value = 6258.500100000001
$ws.Range("M122").Value = -27059.8
$ws.Range("N122").Value = -11158.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5883.9
$ws.Range("I136").Value = 4315.3706
$ws.Range("K136").Value = 12946.1118
$ws.Range("M136").Value = -10396.1118

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 10001.333
$ws.Range("J14").Value = 10002.5
$ws.Range("L14").Value = 10002.5
$ws.Range("N14").Value = -10338.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 30306462
$ws.Range("I81").Value = 45455296
$ws.Range("K81").Value = 90910592
$ws.Range("M81").Value = -90909531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 30306462
$ws.Range("I84").Value = 45455296
$ws.Range("K84").Value = 454552960
$ws.Range("M84").Value = -454547656

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 109999.25
$ws.Range("J118").Value = 109999.25
$ws.Range("L118").Value = 109999.25
$ws.Range("N118").Value = -113313.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2673.132
$ws.Range("I132").Value = 2503.46
$ws.Range("J132").Value = 5501
$ws.Range("K132").Value = 7510.38
$ws.Range("L132").Value = 16503
$ws.Range("M132").Value = -4980.38
$ws.Range("N132").Value = -21563

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1481.093
$ws.Range("I136").Value = 1017.64105
$ws.Range("J136").Value = 5999.75
$ws.Range("K136").Value = 3052.92315
$ws.Range("L136").Value = 17999.25
$ws.Range("M136").Value = -502.9231499999996
$ws.Range("N136").Value = -23099.25
